# edit.ps1 - applies the "strike out finished/obsolete TODO items, add new
# enemy sub-items" revision described in the commit message / xml diff.

function Find-ParaIndex {
    param($doc, $text)
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $pt = $doc.Paragraphs($i).Range.Text.TrimEnd()
        if ($pt -eq $text) {
            return $i
        }
    }
    return -1
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Strike through paragraphs that are now "done"/cut (whole paragraph,
#    every run + the paragraph mark itself gets <w:strike/>).
# ---------------------------------------------------------------------

$strikeTexts = @(
    "Big-ass messaging system for more important messages (scoring, loot)",
    "Layered scrolling, and implicit App2D functions",
    "Game-saving system",
    "A level must be cleared before it can be saved (menu option grayed out?)",
    "Grunt",
    "Ticker (random bullet within range, supercrap)"
)

foreach ($t in $strikeTexts) {
    $i = Find-ParaIndex $d $t
    if ($i -gt 0) {
        $d.Paragraphs($i).Range.Font.StrikeThrough = 1
    }
}

# ---------------------------------------------------------------------
# 2. "Aggro range on enemies, both near (back) and far" -> shortened text,
#    struck through, and split into two new sub-bullets.
# ---------------------------------------------------------------------

$i = Find-ParaIndex $d "Aggro range on enemies, both near (back) and far"
$p = $d.Paragraphs($i)
$p.Range.Text = "Aggro range on enemies"
$p.Range.Font.StrikeThrough = 1

# insert two new ilvl-1 sub-bullets right after it
$p.Range.InsertParagraphAfter()
$pFar = $d.Paragraphs($i + 1)
$pFar.Range.Text = "Far, activate"
$pFar.Range.ListFormat.ListLevelNumber = 2
$pFar.Range.Font.StrikeThrough = 1

$pFar.Range.InsertParagraphAfter()
$pNear = $d.Paragraphs($i + 2)
$pNear.Range.Text = "Near, retract"
$pNear.Range.ListFormat.ListLevelNumber = 2
$pNear.Range.Font.StrikeThrough = 1

# ---------------------------------------------------------------------
# 3. New sub-bullet "Fix main menu with respect to the saving system"
#    right after "A level must be cleared ... out?)" (not struck through).
# ---------------------------------------------------------------------

$i = Find-ParaIndex $d "A level must be cleared before it can be saved (menu option grayed out?)"
$p = $d.Paragraphs($i)
$p.Range.InsertParagraphAfter()
$pFix = $d.Paragraphs($i + 1)
$pFix.Range.Text = "Fix main menu with respect to the saving system"
$pFix.Range.ListFormat.ListLevelNumber = 2

# ---------------------------------------------------------------------
# 4. New top-level bullet "Stick generic enemy properties in their own
#    tag" (struck through) right before "More enemies! E.g".
# ---------------------------------------------------------------------

$i = Find-ParaIndex $d "More enemies! E.g"
$p = $d.Paragraphs($i)
$p.Range.InsertParagraphBefore()
$pStick = $d.Paragraphs($i)
$pStick.Range.Text = "Stick generic enemy properties in their own tag"
$pStick.Range.Font.StrikeThrough = 1

# ---------------------------------------------------------------------
# 5. Mark a lastRenderedPageBreak next to "Impossibler" (cosmetic,
#    matches the diff) and append a new struck-through "Bomber" bullet
#    (ilvl 1) as the very last paragraph of the document.
# ---------------------------------------------------------------------

$last = $d.Paragraphs($d.Paragraphs.Count)
$rStart = $d.Range($last.Range.Start, $last.Range.Start + 11)
$rStart.Collapse(1)
$rStart.InsertXML("<?xml version='1.0'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii='Adobe Caslon Pro' w:hAnsi='Adobe Caslon Pro'/></w:rPr><w:lastRenderedPageBreak/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$pBomber = $d.Paragraphs($d.Paragraphs.Count)
$pBomber.Range.Text = "Bomber"
$pBomber.Range.ListFormat.ListLevelNumber = 2
$pBomber.Range.Font.StrikeThrough = 1

Write-Host "done"
